$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price (D) and 1h volume change (E) values.
# Column D price values that look like plain decimal numbers are entered
# with a leading apostrophe so Excel keeps them as text (matching the
# original inline-string cell contents) instead of auto-converting them
# to numeric values.
$ws.Range("D2").Value = "34.493.57"
$ws.Range("D3").Value = "1.811.84"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'225.71"
$ws.Range("D6").Value = "'0.595"
$ws.Range("E6").Value = "  +2.75%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "'38.32"
$ws.Range("E8").Value = "  +6.77%  "
$ws.Range("D9").Value = "'0.290"
$ws.Range("E9").Value = "  -4.05%  "
$ws.Range("D10").Value = "'0.0676"
$ws.Range("E10").Value = "  -2.69%  "
$ws.Range("E11").Value = "  +0.89%  "
$ws.Range("D12").Value = "2.073.42"
$ws.Range("E12").Value = "  +0.30%  "
$ws.Range("D13").Value = "'11.19"
$ws.Range("E13").Value = "  -2.30%  "
$ws.Range("D14").Value = "1.826.53"
$ws.Range("E14").Value = "  +0.99%  "
$ws.Range("D15").Value = "'0.633"
$ws.Range("E15").Value = "  -1.92%  "
$ws.Range("D16").Value = "34.467.83"
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").Value = "'4.42"
$ws.Range("E17").Value = "  -2.17%  "
$ws.Range("D18").Value = "'68.39"
$ws.Range("E18").Value = "  -1.15%  "
$ws.Range("D19").Value = "'242.46"
$ws.Range("E19").Value = "  -1.40%  "
$ws.Range("E20").Value = "  -2.79%  "
$ws.Range("D21").Value = "'11.24"
$ws.Range("E21").Value = "  -2.15%  "
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").Value = "'4.12"
$ws.Range("E23").Value = "  -1.73%  "
$ws.Range("D24").Value = "'2.22"
$ws.Range("E24").Value = "  +3.77%  "
$ws.Range("D25").Value = "'170.11"
$ws.Range("E25").Value = "  -0.48%  "
$ws.Range("D26").Value = "'7.82"
$ws.Range("E26").Value = "  -1.08%  "
$ws.Range("D27").Value = "'17.60"
$ws.Range("E27").Value = "  +3.56%  "
$ws.Range("E28").Value = "  +1.52%  "
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").Value = "'3.79"
$ws.Range("E30").Value = "  -1.68%  "
$ws.Range("E31").Value = "  -1.55%  "
$ws.Range("D32").Value = "'0.0518"
$ws.Range("E32").Value = "  -2.61%  "
$ws.Range("D33").Value = "'3.84"
$ws.Range("E33").Value = "  -5.60%  "
$ws.Range("D34").Value = "'1.83"
$ws.Range("E34").Value = "  -0.40%  "
$ws.Range("D35").Value = "1.364.14"
$ws.Range("E35").Value = "  -2.46%  "
$ws.Range("D36").Value = "'0.647"
$ws.Range("E36").Value = "  -4.03%  "
$ws.Range("E37").Value = "  -0.71%  "
$ws.Range("D38").Value = "'2.35"
$ws.Range("E38").Value = "  -5.72%  "
$ws.Range("E39").Value = "  -1.82%  "
$ws.Range("E40").Value = "  -1.23%  "
$ws.Range("E41").Value = "  +1.26%  "
$ws.Range("D42").Value = "'0.952"
$ws.Range("E42").Value = "  -1.57%  "
$ws.Range("D43").Value = "'81.70"
$ws.Range("E43").Value = "  -1.28%  "
$ws.Range("E44").Value = "  -0.85%  "
$ws.Range("E45").Value = "  +2.56%  "
$ws.Range("D46").Value = "'0.0510"
$ws.Range("E46").Value = "  +1.43%  "
$ws.Range("D47").Value = "1.974.90"
$ws.Range("E47").Value = "  +0.34%  "
$ws.Range("D48").Value = "'5.77"
$ws.Range("E48").Value = "  -4.74%  "
$ws.Range("E49").Value = "  -0.19%  "
$ws.Range("D50").Value = "'102.04"
$ws.Range("E50").Value = "  -3.10%  "
$ws.Range("E51").Value = "  -5.14%  "
